$wb = $excel.ActiveWorkbook

# Update the "Status" text for the 502a1512-... row (row 3) from
# "Ready for handoff" to "Handback transform failed" everywhere it is
# shown: the Overview summary sheet (zh-cn and de-de status columns) and
# each locale detail sheet's Status column.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Record the handback-vs-handoff filename mismatch error detail in the
# "Error Detail" column (K) for row 3 on both locale sheets.
$wsZhCn.Range("K3").Value = "Handback file name: zpftfklz.cir is different with handoff file name: 502a1512-7933-42f2-b16c-0afa789d4f85.f8fd5796799c9830dc69428de4e2923d5bb53c6f.zh-cn."

$wsDeDe.Range("K3").Value = "Handback file name: zpftfklz.cir is different with handoff file name: 502a1512-7933-42f2-b16c-0afa789d4f85.f8fd5796799c9830dc69428de4e2923d5bb53c6f.de-de."
